$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,7).Value = 0.169654
$ws.Cells.Item(2,8).Value = 0.508962
$ws.Cells.Item(2,9).Value = 0.006094264463659866
$ws.Cells.Item(2,10).Value = 0.006534681579452628
$ws.Cells.Item(2,13).Value = 0.169654
$ws.Cells.Item(2,14).Value = 0.508962
$ws.Cells.Item(2,15).Value = 0.006094264463659866
$ws.Cells.Item(2,16).Value = 0.006534681579452628
$ws.Cells.Item(2,17).Value = 0.028782479716
$ws.Cells.Item(2,18).Value = 0.259042317444
$ws.Cells.Item(2,19).Value = 0.00003714005935302747
$ws.Cells.Item(2,20).Value = 0.00004270206334483749

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,7).Value = 0.169654
$ws.Cells.Item(3,8).Value = 0.508962
$ws.Cells.Item(3,9).Value = 0.006094264463659866
$ws.Cells.Item(3,10).Value = 0.006534681579452628
$ws.Cells.Item(3,13).Value = 21.99231
$ws.Cells.Item(3,14).Value = 65.97693
$ws.Cells.Item(3,15).Value = 0.7900017288527916
$ws.Cells.Item(3,16).Value = 0.8470931604713817
$ws.Cells.Item(3,17).Value = 3.73108336074
$ws.Cells.Item(3,18).Value = 33.57975024666
$ws.Cells.Item(3,19).Value = 0.004814479462377425
$ws.Cells.Item(3,20).Value = 0.005535484071812647

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4,7).Value = 0.169654
$ws.Cells.Item(4,8).Value = 0.508962
$ws.Cells.Item(4,9).Value = 0.006094264463659866
$ws.Cells.Item(4,10).Value = 0.006534681579452628
$ws.Cells.Item(4,13).Value = 0.04769766666666667
$ws.Cells.Item(4,14).Value = 0.143093
$ws.Cells.Item(4,15).Value = 0.001713382501834088
$ws.Cells.Item(4,16).Value = 0.001837204332049573
$ws.Cells.Item(4,17).Value = 0.008092099940666668
$ws.Cells.Item(4,18).Value = 0.072828899466
$ws.Cells.Item(4,19).Value = 0.000010441806093584117
$ws.Cells.Item(4,20).Value = 0.000012005545306334914

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5,7).Value = 0.169654
$ws.Cells.Item(5,8).Value = 0.508962
$ws.Cells.Item(5,9).Value = 0.006094264463659866
$ws.Cells.Item(5,10).Value = 0.006534681579452628
$ws.Cells.Item(5,13).Value = 5.6286445
$ws.Cells.Item(5,14).Value = 11.257289
$ws.Cells.Item(5,15).Value = 0.2021906241817143
$ws.Cells.Item(5,16).Value = 0.1445349536171162
$ws.Cells.Item(5,17).Value = 0.954922054003
$ws.Cells.Item(5,18).Value = 5.729532324018001
$ws.Cells.Item(5,19).Value = 0.0012322031358358286
$ws.Cells.Item(5,20).Value = 0.0009444898989888092

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,7).Value = 21.99231
$ws.Cells.Item(6,8).Value = 65.97693
$ws.Cells.Item(6,9).Value = 0.7900017288527916
$ws.Cells.Item(6,10).Value = 0.8470931604713817
$ws.Cells.Item(6,13).Value = 0.169654
$ws.Cells.Item(6,14).Value = 0.508962
$ws.Cells.Item(6,15).Value = 0.006094264463659866
$ws.Cells.Item(6,16).Value = 0.006534681579452628
$ws.Cells.Item(6,17).Value = 3.73108336074
$ws.Cells.Item(6,18).Value = 33.57975024666
$ws.Cells.Item(6,19).Value = 0.004814479462377425
$ws.Cells.Item(6,20).Value = 0.005535484071812647

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,7).Value = 21.99231
$ws.Cells.Item(7,8).Value = 65.97693
$ws.Cells.Item(7,9).Value = 0.7900017288527916
$ws.Cells.Item(7,10).Value = 0.8470931604713817
$ws.Cells.Item(7,13).Value = 21.99231
$ws.Cells.Item(7,14).Value = 65.97693
$ws.Cells.Item(7,15).Value = 0.7900017288527916
$ws.Cells.Item(7,16).Value = 0.8470931604713817
$ws.Cells.Item(7,17).Value = 483.66169913609997
$ws.Cells.Item(7,18).Value = 4352.9552922248995
$ws.Cells.Item(7,19).Value = 0.6241027315903998
$ws.Cells.Item(7,20).Value = 0.7175668225173941

# Row 8: FAPs -> Inflammatory-Mac
$ws.Cells.Item(8,7).Value = 21.99231
$ws.Cells.Item(8,8).Value = 65.97693
$ws.Cells.Item(8,9).Value = 0.7900017288527916
$ws.Cells.Item(8,10).Value = 0.8470931604713817
$ws.Cells.Item(8,13).Value = 0.04769766666666667
$ws.Cells.Item(8,14).Value = 0.143093
$ws.Cells.Item(8,15).Value = 0.001713382501834088
$ws.Cells.Item(8,16).Value = 0.001837204332049573
$ws.Cells.Item(8,17).Value = 1.0489818716100001
$ws.Cells.Item(8,18).Value = 9.440836844489999
$ws.Cells.Item(8,19).Value = 0.001353575138635051
$ws.Cells.Item(8,20).Value = 0.0015562832240675866

# Row 9: FAPs -> MuSCs
$ws.Cells.Item(9,7).Value = 21.99231
$ws.Cells.Item(9,8).Value = 65.97693
$ws.Cells.Item(9,9).Value = 0.7900017288527916
$ws.Cells.Item(9,10).Value = 0.8470931604713817
$ws.Cells.Item(9,13).Value = 5.6286445
$ws.Cells.Item(9,14).Value = 11.257289
$ws.Cells.Item(9,15).Value = 0.2021906241817143
$ws.Cells.Item(9,16).Value = 0.1445349536171162
$ws.Cells.Item(9,17).Value = 123.786894723795
$ws.Cells.Item(9,18).Value = 742.72136834277
$ws.Cells.Item(9,19).Value = 0.15973094266137935
$ws.Cells.Item(9,20).Value = 0.12243457065810753

# Row 10: Inflammatory-Mac -> ECs
$ws.Cells.Item(10,7).Value = 0.04769766666666667
$ws.Cells.Item(10,8).Value = 0.143093
$ws.Cells.Item(10,9).Value = 0.001713382501834088
$ws.Cells.Item(10,10).Value = 0.001837204332049573
$ws.Cells.Item(10,13).Value = 0.169654
$ws.Cells.Item(10,14).Value = 0.508962
$ws.Cells.Item(10,15).Value = 0.006094264463659866
$ws.Cells.Item(10,16).Value = 0.006534681579452628
$ws.Cells.Item(10,17).Value = 0.008092099940666668
$ws.Cells.Item(10,18).Value = 0.072828899466
$ws.Cells.Item(10,19).Value = 0.000010441806093584117
$ws.Cells.Item(10,20).Value = 0.000012005545306334914

# Row 11: Inflammatory-Mac -> FAPs
$ws.Cells.Item(11,7).Value = 0.04769766666666667
$ws.Cells.Item(11,8).Value = 0.143093
$ws.Cells.Item(11,9).Value = 0.001713382501834088
$ws.Cells.Item(11,10).Value = 0.001837204332049573
$ws.Cells.Item(11,13).Value = 21.99231
$ws.Cells.Item(11,14).Value = 65.97693
$ws.Cells.Item(11,15).Value = 0.7900017288527916
$ws.Cells.Item(11,16).Value = 0.8470931604713817
$ws.Cells.Item(11,17).Value = 1.0489818716100001
$ws.Cells.Item(11,18).Value = 9.440836844489999
$ws.Cells.Item(11,19).Value = 0.001353575138635051
$ws.Cells.Item(11,20).Value = 0.0015562832240675866

# Row 12: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(12,7).Value = 0.04769766666666667
$ws.Cells.Item(12,8).Value = 0.143093
$ws.Cells.Item(12,9).Value = 0.001713382501834088
$ws.Cells.Item(12,10).Value = 0.001837204332049573
$ws.Cells.Item(12,13).Value = 0.04769766666666667
$ws.Cells.Item(12,14).Value = 0.143093
$ws.Cells.Item(12,15).Value = 0.001713382501834088
$ws.Cells.Item(12,16).Value = 0.001837204332049573
$ws.Cells.Item(12,17).Value = 0.002275067405444445
$ws.Cells.Item(12,18).Value = 0.020475606649
$ws.Cells.Item(12,19).Value = 0.0000029356795975912384
$ws.Cells.Item(12,20).Value = 0.0000033753197577017178

# Row 13: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(13,7).Value = 0.04769766666666667
$ws.Cells.Item(13,8).Value = 0.143093
$ws.Cells.Item(13,9).Value = 0.001713382501834088
$ws.Cells.Item(13,10).Value = 0.001837204332049573
$ws.Cells.Item(13,13).Value = 5.6286445
$ws.Cells.Item(13,14).Value = 11.257289
$ws.Cells.Item(13,15).Value = 0.2021906241817143
$ws.Cells.Item(13,16).Value = 0.1445349536171162
$ws.Cells.Item(13,17).Value = 0.2684732091461667
$ws.Cells.Item(13,18).Value = 1.610839254877
$ws.Cells.Item(13,19).Value = 0.0003464298775078615
$ws.Cells.Item(13,20).Value = 0.00026554024291795

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14,7).Value = 5.6286445
$ws.Cells.Item(14,8).Value = 11.257289
$ws.Cells.Item(14,9).Value = 0.2021906241817143
$ws.Cells.Item(14,10).Value = 0.1445349536171162
$ws.Cells.Item(14,13).Value = 0.169654
$ws.Cells.Item(14,14).Value = 0.508962
$ws.Cells.Item(14,15).Value = 0.006094264463659866
$ws.Cells.Item(14,16).Value = 0.006534681579452628
$ws.Cells.Item(14,17).Value = 0.954922054003
$ws.Cells.Item(14,18).Value = 5.729532324018001
$ws.Cells.Item(14,19).Value = 0.0012322031358358286
$ws.Cells.Item(14,20).Value = 0.0009444898989888092

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15,7).Value = 5.6286445
$ws.Cells.Item(15,8).Value = 11.257289
$ws.Cells.Item(15,9).Value = 0.2021906241817143
$ws.Cells.Item(15,10).Value = 0.1445349536171162
$ws.Cells.Item(15,13).Value = 21.99231
$ws.Cells.Item(15,14).Value = 65.97693
$ws.Cells.Item(15,15).Value = 0.7900017288527916
$ws.Cells.Item(15,16).Value = 0.8470931604713817
$ws.Cells.Item(15,17).Value = 123.786894723795
$ws.Cells.Item(15,18).Value = 742.72136834277
$ws.Cells.Item(15,19).Value = 0.15973094266137935
$ws.Cells.Item(15,20).Value = 0.12243457065810753

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16,7).Value = 5.6286445
$ws.Cells.Item(16,8).Value = 11.257289
$ws.Cells.Item(16,9).Value = 0.2021906241817143
$ws.Cells.Item(16,10).Value = 0.1445349536171162
$ws.Cells.Item(16,13).Value = 0.04769766666666667
$ws.Cells.Item(16,14).Value = 0.143093
$ws.Cells.Item(16,15).Value = 0.001713382501834088
$ws.Cells.Item(16,16).Value = 0.001837204332049573
$ws.Cells.Item(16,17).Value = 0.2684732091461667
$ws.Cells.Item(16,18).Value = 1.610839254877
$ws.Cells.Item(16,19).Value = 0.0003464298775078615
$ws.Cells.Item(16,20).Value = 0.00026554024291795

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17,7).Value = 5.6286445
$ws.Cells.Item(17,8).Value = 11.257289
$ws.Cells.Item(17,9).Value = 0.2021906241817143
$ws.Cells.Item(17,10).Value = 0.1445349536171162
$ws.Cells.Item(17,13).Value = 5.6286445
$ws.Cells.Item(17,14).Value = 11.257289
$ws.Cells.Item(17,15).Value = 0.2021906241817143
$ws.Cells.Item(17,16).Value = 0.1445349536171162
$ws.Cells.Item(17,17).Value = 31.68163890738025
$ws.Cells.Item(17,18).Value = 126.726555629521
$ws.Cells.Item(17,19).Value = 0.04088104850699123
$ws.Cells.Item(17,20).Value = 0.020890352817101936
